$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("model_1_9_6", 0.9429760926646541, 0.767918987163799, 0.6552020937782022, 0.7750383533850937, 0.7159562135324031, 0.3813191309224081, 1.551926801470479, 0.3425257021022688, 0.1818469795453996, 0.2621865092072034, 0.9467346377056539, 0.617510429808605, 1.036988480433738, 0.6437991272283683, 123.9282372809774, 198.2796625979376)
    3 = @("model_1_9_7", 0.943568724768018, 0.767917064191121, 0.6532844458700133, 0.762073732033557, 0.7092476940407448, 0.3773561973183937, 1.551939660397631, 0.3444307127890646, 0.1923268869838552, 0.2683788055758022, 0.9397426380818825, 0.6142932502627664, 1.036604070420745, 0.6404449856887836, 123.9491314342457, 198.300556751206)
    4 = @("model_1_9_5", 0.9423340397840376, 0.7678729681171016, 0.6563896345728446, 0.7872823818445287, 0.7219565169460784, 0.3856125414914719, 1.552234531047643, 0.3413459871529606, 0.1719495609128467, 0.2566479314203766, 0.9523670069824894, 0.6209770861243367, 1.037404947167111, 0.6474133662804286, 123.9058443843197, 198.2572697012799)
    5 = @("model_1_9_8", 0.9441145449408084, 0.7678692029519372, 0.6506690121495615, 0.7484570114997497, 0.7018779882641342, 0.3737062953096736, 1.552259708724435, 0.3470289109081198, 0.2033339165715493, 0.2751814096935801, 0.932028931607931, 0.611315217632993, 1.036250024903259, 0.6373401720445834, 123.9685701934821, 198.3199955104423)
    6 = @("model_1_9_9", 0.9446156621354717, 0.7677775121951159, 0.6473831375496775, 0.7342558746310599, 0.6938915349761421, 0.3703553223215279, 1.552872845237495, 0.3502931317285946, 0.2148133571097079, 0.2825533023674722, 0.9355976335984386, 0.6085682560909071, 1.035924975912127, 0.6344762666626566, 123.9865848071351, 198.3380101240954)
    7 = @("model_1_9_4", 0.9416400155044878, 0.7677768071900872, 0.6568101047983755, 0.7987301746452293, 0.7271956609160837, 0.3902534850445097, 1.552877559609368, 0.3409282878090428, 0.1626957766585354, 0.2518119415689092, 0.9565370826967162, 0.6247027173340209, 1.03785512507817, 0.6512976053238289, 123.8819175786789, 198.2333428956391)
    8 = @("model_1_9_10", 0.945074301774531, 0.767643947102157, 0.6434552075578657, 0.7195303900762435, 0.6853301212934974, 0.3672883969432967, 1.553765995629123, 0.3541951768221877, 0.2267166523110436, 0.2904559120152543, 0.937924939501726, 0.6060432302594401, 1.035627479930034, 0.63184374525402, 124.0032158335816, 198.3546411505419)
    9 = @("model_1_9_3", 0.9408912658538993, 0.7676283781789921, 0.6564275349653981, 0.8093017378154086, 0.731618479596379, 0.3952603774056703, 1.553870105090014, 0.3413083365224461, 0.1541502896366462, 0.2477294604660341, 0.9591358907459677, 0.6286973655151342, 1.038340800527201, 0.6554623139480739, 123.8564210958828, 198.2078464128431)
    10 = @("model_1_9_11", 0.9454924826743007, 0.76747046670894, 0.6389082522844043, 0.7043357906259453, 0.6762302804029987, 0.3644920193409902, 1.554926059817368, 0.3587121678461068, 0.2389991549376751, 0.2988555166927289, 0.9390963908592554, 0.6037317445198572, 1.035356227454508, 0.629433854781068, 124.0185012451704, 198.3699265621307)
    11 = @("model_1_9_2", 0.9400848771624585, 0.767425540561252, 0.655199959793459, 0.8189203509645643, 0.7351696155321272, 0.4006527022984132, 1.555226481173824, 0.3425278220241397, 0.1463751165130968, 0.2444515857894117, 0.9600561645142789, 0.6329713281803633, 1.038863863462189, 0.6599182280522842, 123.8293206119901, 198.1807459289503)
    12 = @("model_1_9_12", 0.9458719403464211, 0.7672588703560209, 0.6337667809640263, 0.6887236584750815, 0.6666273859541234, 0.3619545841402793, 1.55634100560365, 0.363819757080475, 0.2516191686983232, 0.3077194647661823, 0.9391902801646153, 0.6016266152193396, 1.035110092748267, 0.6272391057680088, 124.0324730664059, 198.3838983833661)
    13 = @("model_1_9_1", 0.9392180097524032, 0.7671661612992476, 0.6530838341850358, 0.8275026938623274, 0.7377890516648509, 0.4064494486610103, 1.556960951493179, 0.3446300024512317, 0.1394376089118446, 0.2420337162620727, 0.9591736734753081, 0.6375338804024538, 1.039426155836279, 0.6646750175682542, 123.8005914305296, 198.1520167474899)
    14 = @("model_1_9_13", 0.9462144593359327, 0.7670109266047431, 0.6280491814515806, 0.6727392398929424, 0.6565516182800506, 0.3596641580802603, 1.55799900661012, 0.3694996778456557, 0.2645401189255731, 0.3170198982905371, 0.9382925675938981, 0.5997200664312144, 1.034887918268584, 0.6252513912508731, 124.0451691543402, 198.3965944713005)
    15 = @("model_1_9_0", 0.9382877274148292, 0.7668481426425062, 0.6500346625698117, 0.834974940194699, 0.7394200085638106, 0.4126702509359239, 1.55908754371515, 0.3476590801500478, 0.1333974440819392, 0.2405282621159935, 0.9563792143191479, 0.6423941554341258, 1.040029582217408, 0.6697422045701809, 123.7702128576304, 198.1216381745906)
    16 = @("model_1_9_14", 0.94652148702953, 0.7667283678621292, 0.621775064204954, 0.6564215357488965, 0.646030444893787, 0.3576110624794427, 1.559888478223043, 0.3757324489158789, 0.2777304793997285, 0.3267314633883459, 0.9701609367167542, 0.5980059050539909, 1.034688765170034, 0.623464254475018, 124.0566186016321, 198.4080439185924)
    17 = @("model_1_9_15", 0.9467945682409398, 0.7664127411207148, 0.6149642117186374, 0.6398104227988525, 0.6350919680028524, 0.3557849671613221, 1.561999075696209, 0.3824984181624843, 0.2911580158811975, 0.336828220327632, 1.006161483577694, 0.5964771304596029, 1.034511631411282, 0.6218703967811472, 124.0668575109087, 198.4182828278689)
    18 = @("model_1_9_16", 0.9470349935344354, 0.7660656009332028, 0.6076313472000165, 0.6229356579763414, 0.6237571475568331, 0.3541772421166601, 1.564320403728535, 0.3897829594034213, 0.3047986744542809, 0.3472908220348911, 1.042103833150674, 0.5951279208007806, 1.034355679869555, 0.6204637484738997, 124.0759156144068, 198.427340931367)
    19 = @("model_1_9_17", 0.9472439176602973, 0.7656883283081996, 0.5997917091499925, 0.6058274283471092, 0.6120474027824947, 0.3527801655250148, 1.566843227509112, 0.3975709345588874, 0.3186280535073718, 0.3580994975009933, 1.077960856127307, 0.5939529994242093, 1.034220161517645, 0.6192388082619054, 124.0838203531942, 198.4352456701544)
    20 = @("model_1_9_18", 0.9474225570555646, 0.7652823334745069, 0.5914616861108251, 0.5885103396628371, 0.5999820986338283, 0.351585603066304, 1.569558116831449, 0.4058460630864474, 0.3326262630650451, 0.3692363719125126, 1.113710587031129, 0.5929465431776325, 1.034104287315309, 0.6181895050892601, 124.0906041218397, 198.4420294388)
    21 = @("model_1_9_19", 0.947572000834007, 0.7648488585285822, 0.5826524630205365, 0.5710104205596978, 0.587579321650165, 0.3505862718317334, 1.572456765790816, 0.4145972337563633, 0.346772262968052, 0.3806847504961895, 1.14933280511743, 0.5921032611223598, 1.034007350810374, 0.61731032277106, 124.0962969269492, 198.4477222439094)
    22 = @("model_1_9_20", 0.9476931912938307, 0.7643891587604594, 0.573379337744433, 0.5533445140707016, 0.5748540744138508, 0.3497758706688618, 1.575530780256976, 0.4238092495156366, 0.3610524382081213, 0.3924307851240881, 1.184812054868087, 0.5914185241171109, 1.03392874078238, 0.6165964350938967, 124.1009253982807, 198.452350715241)
    23 = @("model_1_9_21", 0.9477870943532246, 0.763904383285781, 0.5636553831581664, 0.5355332309613781, 0.5618226200806071, 0.3491479404783112, 1.578772476088336, 0.433469123638332, 0.3754501281432547, 0.4044594640964972, 1.220136622136913, 0.590887417769503, 1.0338678306898, 0.6160427184495294, 124.104519096761, 198.4559444137213)
    24 = @("model_1_9_22", 0.9478546389526019, 0.7633956265199686, 0.553492283133669, 0.5175903381236805, 0.5484969324430684, 0.3486962694312604, 1.582174534923934, 0.4435652492487563, 0.3899542043533099, 0.4167597350086738, 1.255283887871513, 0.5905050968715345, 1.033824017976691, 0.615644121359421, 124.1071080477258, 198.458533364686)
    25 = @("model_1_9_23", 0.9478965783402495, 0.7628639014107368, 0.5429020376943048, 0.499532390774913, 0.5348911968027762, 0.3484158205529482, 1.585730183177735, 0.4540857054032092, 0.4045512844848599, 0.4293185041234179, 1.290250280072065, 0.5902675838574809, 1.033796814049568, 0.6153964969246393, 124.1087172517403, 198.4601425687005)
    26 = @("model_1_9_24", 0.9479138015513239, 0.7623102162283772, 0.5318952576226522, 0.481374807055816, 0.5210176051431754, 0.3483006488611881, 1.58943267854166, 0.465019951243736, 0.419228905336439, 0.4421245176350419, 1.325025562188786, 0.5901700169113882, 1.033785642236979, 0.61529477635166, 124.1093784773485, 198.4608037943088)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item([int]$r, 1).Value = $row[0]
    for ($c = 0; $c -lt 16; $c++) {
        $ws.Cells.Item([int]$r, $c + 2).Value = $row[$c + 1]
    }
}
